$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 5858510
$ws.Range("I98").Value = 16925.908
$ws.Range("J98").Value = 13890688
$ws.Range("K98").Value = 16925.908
$ws.Range("L98").Value = 13890688
$ws.Range("M98").Value = -15427.908
$ws.Range("N98").Value = -13893684

$ws.Range("H122").Value = 5858510
$ws.Range("I122").Value = 16925.908
$ws.Range("J122").Value = 13890688
$ws.Range("K122").Value = 50777.724
$ws.Range("L122").Value = 41672064
$ws.Range("M122").Value = -48327.724
$ws.Range("N122").Value = -41676964

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H125").Value = 1337.3334
$ws.Range("I125").Value = 516
$ws.Range("J125").Value = 2980
$ws.Range("K125").Value = 4644
$ws.Range("L125").Value = 26820
$ws.Range("M125").Value = -2184
$ws.Range("N125").Value = -31740

$ws.Range("H126").Value = 48660
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 48660
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 48660
$ws.Range("N126").Value = -58540

$ws.Range("H127").Value = 1019.25
$ws.Range("I127").Value = 359
$ws.Range("J127").Value = 3000
$ws.Range("K127").Value = 1077
$ws.Range("L127").Value = 9000
$ws.Range("M127").Value = 3883
$ws.Range("N127").Value = -18920

$ws.Range("H129").Value = 972.7143
$ws.Range("I129").Value = 410
$ws.Range("J129").Value = 1066.5
$ws.Range("K129").Value = 1230
$ws.Range("L129").Value = 3199.5
$ws.Range("M129").Value = 3770
$ws.Range("N129").Value = -13199.5

$ws.Range("H130").Value = 37134.285
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 37134.285
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 37134.285
$ws.Range("N130").Value = -47174.285

$ws.Range("H131").Value = 5707.5
$ws.Range("I131").Value = 2332
$ws.Range("J131").Value = 11333.333
$ws.Range("K131").Value = 6996
$ws.Range("L131").Value = 33999.999
$ws.Range("M131").Value = -1956
$ws.Range("N131").Value = -44079.999

$ws.Range("H132").Value = 3269.3333
$ws.Range("I132").Value = 3184.6191
$ws.Range("J132").Value = 3862.3333
$ws.Range("K132").Value = 9553.8573
$ws.Range("L132").Value = 11586.9999
$ws.Range("M132").Value = -7023.8573
$ws.Range("N132").Value = -16646.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1658.7273
$ws.Range("I110").Value = 1310
$ws.Range("J110").Value = 3228
$ws.Range("K110").Value = 1310
$ws.Range("L110").Value = 3228
$ws.Range("M110").Value = 735
$ws.Range("N110").Value = -7318

$ws.Range("H132").Value = 2460.8333
$ws.Range("I132").Value = 1430.25
$ws.Range("J132").Value = 3638.6428
$ws.Range("K132").Value = 4290.75
$ws.Range("L132").Value = 10915.9284
$ws.Range("M132").Value = -1760.75
$ws.Range("N132").Value = -15975.9284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2524.0645
$ws.Range("I31").Value = 2093.3044
$ws.Range("J31").Value = 3762.5
$ws.Range("K31").Value = 2093.3044
$ws.Range("L31").Value = 3762.5
$ws.Range("M31").Value = -1798.3044
$ws.Range("N31").Value = -4352.5

$ws.Range("H34").Value = 2524.0645
$ws.Range("I34").Value = 2093.3044
$ws.Range("J34").Value = 3762.5
$ws.Range("K34").Value = 2093.3044
$ws.Range("L34").Value = 3762.5
$ws.Range("M34").Value = -1891.3044
$ws.Range("N34").Value = -4166.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 6015.5
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 6015.5
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 6015.5
$ws.Range("N47").Value = -7151.5

$ws.Range("H132").Value = 4740.175
$ws.Range("I132").Value = 5814.885
$ws.Range("J132").Value = 2744.2856
$ws.Range("K132").Value = 17444.655
$ws.Range("L132").Value = 8232.856800000001
$ws.Range("M132").Value = -14914.655
$ws.Range("N132").Value = -13292.8568

$ws.Range("H133").Value = 58900
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 58900
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 58900
$ws.Range("N133").Value = -69020

$ws.Range("H135").Value = 47700
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 47700
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 47700
$ws.Range("N135").Value = -57840

$ws.Range("H139").Value = 25000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 25000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 25000
$ws.Range("N139").Value = -35280

$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 39582.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 39582.25
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 39582.25
$ws.Range("N141").Value = -49942.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2357.7778
$ws.Range("I7").Value = 2363.889
$ws.Range("J7").Value = 2351.6667
$ws.Range("K7").Value = 2363.889
$ws.Range("L7").Value = 2351.6667
$ws.Range("M7").Value = -2251.889
$ws.Range("N7").Value = -2575.6667

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H22").Value = 2568.6924
$ws.Range("I22").Value = 3440.6177
$ws.Range("J22").Value = 921.7222
$ws.Range("K22").Value = 3440.6177
$ws.Range("L22").Value = 921.7222
$ws.Range("M22").Value = -3145.6177
$ws.Range("N22").Value = -1511.7222

$ws.Range("H27").Value = 2568.6924
$ws.Range("I27").Value = 3440.6177
$ws.Range("J27").Value = 921.7222
$ws.Range("K27").Value = 3440.6177
$ws.Range("L27").Value = 921.7222
$ws.Range("M27").Value = -3333.6177
$ws.Range("N27").Value = -1135.7222

$ws.Range("H122").Value = 3284.4
$ws.Range("I122").Value = 2963.75
$ws.Range("J122").Value = 4567
$ws.Range("K122").Value = 8891.25
$ws.Range("L122").Value = 13701
$ws.Range("M122").Value = -6441.25
$ws.Range("N122").Value = -18601

$ws.Range("H123").Value = 54980
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 54980
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 54980
$ws.Range("N123").Value = -64780

$ws.Range("H124").Value = 28000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 28000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 28000
$ws.Range("N124").Value = -37820

$ws.Range("H125").Value = 31357.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 31357.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 31357.5
$ws.Range("N125").Value = -41197.5

$ws.Range("H126").Value = 2357.7778
$ws.Range("I126").Value = 2363.889
$ws.Range("J126").Value = 2351.6667
$ws.Range("K126").Value = 7091.667
$ws.Range("L126").Value = 7055.000100000001
$ws.Range("M126").Value = -4621.667
$ws.Range("N126").Value = -11995.0001

$ws.Range("H127").Value = 62566.45
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 62566.45
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 62566.45
$ws.Range("N127").Value = -72486.45

$ws.Range("H129").Value = 28000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 28000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 28000
$ws.Range("N129").Value = -38000

$ws.Range("H130").Value = 39990
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 39990
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 39990
$ws.Range("N130").Value = -50030

$ws.Range("H131").Value = 15732.6
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 15732.6
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 15732.6
$ws.Range("N131").Value = -25812.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1740.762
$ws.Range("I122").Value = 1703.4615
$ws.Range("J122").Value = 1801.375
$ws.Range("K122").Value = 5110.3845
$ws.Range("L122").Value = 5404.125
$ws.Range("M122").Value = -2660.3845
$ws.Range("N122").Value = -10304.125

$ws.Range("H123").Value = 30000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 30000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H126").Value = 1737.5667
$ws.Range("I126").Value = 1281.0588
$ws.Range("J126").Value = 2334.5386
$ws.Range("K126").Value = 3843.1764
$ws.Range("L126").Value = 7003.6158
$ws.Range("M126").Value = -1373.1764
$ws.Range("N126").Value = -11943.6158

$ws.Range("H127").Value = 57469.668
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 57469.668
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 57469.668
$ws.Range("N127").Value = -67389.66800000001

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H129").Value = 19714.5
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 19714.5
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 19714.5
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -29714.5

$ws.Range("H131").Value = 49958
$ws.Range("I131").Value = 25000
$ws.Range("J131").Value = 56197.5
$ws.Range("K131").Value = 25000
$ws.Range("L131").Value = 56197.5
$ws.Range("M131").Value = -19960
$ws.Range("N131").Value = -66277.5
